# Feature: add arrows (arrow_n). Fixed bugs, removed unnecessary code.
#
# The "meta" worksheet stores key/value pairs in columns A/B (column A uses
# the bold "key" style, index 1). The sheet used to end with a lone,
# value-less "key style" cell (A5) that acted as a trailing placeholder row.
# This change inserts a new "style" / "default" key/value pair right before
# that placeholder row, pushing the placeholder row down by one (from row 5
# to row 6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("meta")

# Insert a new blank row at row 5 (this shifts the existing placeholder row,
# together with its formatting, from row 5 down to row 6).
$ws.Rows.Item(5).Insert()

# Populate the newly freed row 5 with the new "style" = "default" entry.
$ws.Cells.Item(5, 1).Value = "style"
$ws.Cells.Item(5, 2).Value = "default"
